$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999817531449
$ws.Range("A2").Value = 0.99822233349817202
$ws.Range("A3").Value = 0.99227080263766787
$ws.Range("A4").Value = 0.99441002924156452
$ws.Range("A5").Value = 0.98336561852589133
$ws.Range("A6").Value = 0.95843640201629809
$ws.Range("A7").Value = 0.95345995971104802
$ws.Range("A8").Value = 0.94587886485025763
$ws.Range("A9").Value = 0.93721240015165153
$ws.Range("A10").Value = 0.92945457603075743
$ws.Range("A11").Value = 0.92827547268278288
$ws.Range("A12").Value = 0.9257687968476922
$ws.Range("A13").Value = 0.91644986095550329
$ws.Range("A14").Value = 0.91339104775929836
$ws.Range("A15").Value = 0.91166807040383646
$ws.Range("A16").Value = 0.90916152964660812
$ws.Range("A17").Value = 0.90545351735652657
$ws.Range("A18").Value = 0.90434457682790159
$ws.Range("A19").Value = 0.99716906444227704
$ws.Range("A20").Value = 0.99005186020203428
$ws.Range("A21").Value = 0.98865333733791827
$ws.Range("A22").Value = 0.98738882724044941
$ws.Range("A23").Value = 0.98456851651393085
$ws.Range("A24").Value = 0.97154806979578834
$ws.Range("A25").Value = 0.96509113368195321
$ws.Range("A26").Value = 0.95594219390004898
$ws.Range("A27").Value = 0.95271361494043882
$ws.Range("A28").Value = 0.94085270250095543
$ws.Range("A29").Value = 0.93270000373691619
$ws.Range("A30").Value = 0.92889623821040013
$ws.Range("A31").Value = 0.93083089553270604
$ws.Range("A32").Value = 0.93272006816360853
$ws.Range("A33").Value = 0.93220003983732069
